$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally held 20 data rows (A2:C21, w = 238..4760 in steps of
# 238). The data set was regenerated with a coarser, rounder step (w =
# 500..9500 in steps of 500) and one fewer sample, so the used range shrinks
# from A1:C21 to A1:C20 and every data row's values change.

# Drop the now-unused last row (old row 21) first so the sheet's dimension
# becomes A1:C20.
$ws.Rows.Item(21).Delete()

# New w / deltaF / fit=a*w^b values for rows 2-20.
$newData = @(
    @(500,  3.069653889284591,  3.335148170911226),
    @(1000, 4.429804532502494,  4.518919412959403),
    @(1500, 5.334960779615662,  5.397622850611862),
    @(2000, 6.145941898393136,  6.122856201390911),
    @(2500, 6.816020514158049,  6.751847762248596),
    @(3000, 7.395939851963604,  7.313444990601127),
    @(3500, 7.940500039650714,  7.824559330911549),
    @(4000, 8.373555978848787,  8.29609130789072),
    @(4500, 8.795707876050505,  8.735543653368561),
    @(5000, 9.198734421718591,  9.148335954691776),
    @(5500, 9.609058089643931,  9.538531009233115),
    @(6000, 9.964032169466703,  9.909265450898557),
    @(6500, 10.28937092092352,  10.2630195234783),
    @(7000, 10.61427501517862,  10.60179375738146),
    @(7500, 10.89400910350406,  10.92722898587394),
    @(8000, 11.21266361154406,  11.24069041066572),
    @(8500, 11.48195363740862,  11.54332804302194),
    @(9000, 11.76149241736466,  11.8361211481576),
    @(9500, 12.00139562254284,  12.11991157349273)
)

$row = 2
foreach ($r in $newData) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row++
}
